# Generate Report for Handback
# Update the timestamp text values that record when the handoff / handback
# xliff files were generated, as part of regenerating the handback status
# report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 07:08:29"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 07:08:23"
$wsZhCn.Range("K2").Value = "2016-08-30 07:08:40"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 07:08:29"
$wsDeDe.Range("K2").Value = "2016-08-30 07:08:47"
